$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue 'D2' '326.02'
Set-TextValue 'E2' '-1.26%'
Set-TextValue 'G2' '5'
Set-TextValue 'D3' '44.39'
Set-TextValue 'E3' '1.44%'
Set-TextValue 'G3' '5'
Set-TextValue 'D4' '5.534'
Set-TextValue 'E4' '-4.87%'
Set-TextValue 'G4' '5'
Set-TextValue 'G5' '5'
Set-TextValue 'D6' '8.694'
Set-TextValue 'E6' '-0.99%'
Set-TextValue 'G6' '5'
Set-TextValue 'B7' 'GateToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.335'
Set-TextValue 'E7' '-3.67%'
Set-TextValue 'G7' '5'
Set-TextValue 'B8' 'FTXToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D8' '1.895'
Set-TextValue 'E8' '-3.29%'
Set-TextValue 'G8' '5'
Set-TextValue 'G9' '5'
Set-TextValue 'D10' '0.9472'
Set-TextValue 'E10' '1.37%'
Set-TextValue 'G10' '5'
Set-TextValue 'D11' '0.1178'
Set-TextValue 'G11' '5'
Set-TextValue 'D12' '0.1899'
Set-TextValue 'E12' '-2.88%'
Set-TextValue 'G12' '5'
Set-TextValue 'D13' '0.09983'
Set-TextValue 'E13' '6.61%'
Set-TextValue 'G13' '5'
Set-TextValue 'D14' '0.04172'
Set-TextValue 'E14' '4.97%'
Set-TextValue 'G14' '5'
Set-TextValue 'D15' '0.1064'
Set-TextValue 'E15' '0.06%'
Set-TextValue 'G15' '5'
Set-TextValue 'D16' '0.001268'
Set-TextValue 'E16' '-2.82%'
Set-TextValue 'G16' '5'
Set-TextValue 'D17' '0.005909'
Set-TextValue 'E17' '-0.88%'
Set-TextValue 'G17' '5'
Set-TextValue 'E18' '2.52%'
Set-TextValue 'G18' '5'
Set-TextValue 'E19' '-0.66%'
Set-TextValue 'G19' '5'
Set-TextValue 'D20' '8.381'
Set-TextValue 'E20' '-8.00%'
Set-TextValue 'G20' '5'
Set-TextValue 'D21' '0.1372'
Set-TextValue 'E21' '0.12%'
Set-TextValue 'G21' '5'
Set-TextValue 'D22' '0.2662'
Set-TextValue 'E22' '3.57%'
Set-TextValue 'G22' '5'
Set-TextValue 'D23' '0.04249'
Set-TextValue 'E23' '-3.04%'
Set-TextValue 'G23' '5'
Set-TextValue 'D24' '0.001240'
Set-TextValue 'E24' '-1.47%'
Set-TextValue 'G24' '5'
Set-TextValue 'D25' '0.004603'
Set-TextValue 'E25' '3.83%'
Set-TextValue 'G25' '5'
Set-TextValue 'D26' '0.0001235'
Set-TextValue 'E26' '3.66%'
Set-TextValue 'G26' '5'
Set-TextValue 'D27' '0.0003996'
Set-TextValue 'G27' '5'
Set-TextValue 'G28' '5'
Set-TextValue 'G29' '5'
Set-TextValue 'G30' '5'
Set-TextValue 'G31' '5'
Set-TextValue 'G32' '5'
Set-TextValue 'G33' '5'
Set-TextValue 'G34' '5'
Set-TextValue 'G35' '5'
Set-TextValue 'G36' '5'
Set-TextValue 'G37' '5'
Set-TextValue 'G38' '5'
Set-TextValue 'D39' '0.02650'
Set-TextValue 'E39' '-5.71%'
Set-TextValue 'G39' '5'
Set-TextValue 'D40' '0.05535'
Set-TextValue 'E40' '-1.39%'
Set-TextValue 'G40' '5'
Set-TextValue 'E41' '25.01%'
Set-TextValue 'G41' '5'
Set-TextValue 'D42' '0.007676'
Set-TextValue 'E42' '-3.11%'
Set-TextValue 'G42' '5'
Set-TextValue 'D43' '0.1392'
Set-TextValue 'E43' '-2.35%'
Set-TextValue 'G43' '5'
Set-TextValue 'E44' '-4.37%'
Set-TextValue 'G44' '5'
Set-TextValue 'D45' '0.008679'
Set-TextValue 'E45' '-16.63%'
Set-TextValue 'G45' '5'
Set-TextValue 'D46' '0.00007122'
Set-TextValue 'E46' '-1.04%'
Set-TextValue 'G46' '5'
Set-TextValue 'D47' '0.00000000751'
Set-TextValue 'E47' '0.07%'
Set-TextValue 'G47' '5'
Set-TextValue 'D48' '0.003431'
Set-TextValue 'E48' '-13.85%'
Set-TextValue 'G48' '5'
Set-TextValue 'D49' '0.002274'
Set-TextValue 'G49' '5'
Set-TextValue 'D50' '0.00002104'
Set-TextValue 'E50' '0.07%'
Set-TextValue 'G50' '5'
Set-TextValue 'D51' '0.0002004'
Set-TextValue 'E51' '0.07%'
Set-TextValue 'G51' '5'
